# This script refreshes the cached market-price / profit figures (columns H-N)
# on the leve-profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to the
# latest values pulled by the scheduled market-data runner.
#
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
#
# Note: M (LeveProfitNQ) only exists when K (LevePriceNQ) is non-zero, and
# N (LeveProfitHQ) only exists when L (LevePriceHQ) is non-zero - matching the
# source data's convention of omitting profit figures when no price is available.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 118.81481
$ws.Range("I33").Value = 78.13636
$ws.Range("J33").Value = 297.8
$ws.Range("K33").Value = 78.13636
$ws.Range("L33").Value = 297.8
$ws.Range("M33").Value = 150.86364
$ws.Range("N33").Value = -755.8

# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 2860
$ws.Range("I40").Value = 2825
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2825
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2650
$ws.Range("N40").Value = -3350

# Row 105: Ultimate Official Strategy Guide / Gazelleskin Codex
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = ""

# Row 125: Body over Mind / Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 1125
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").Value = ""

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 30309328
$ws.Range("I132").Value = 37043844
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 111131532
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -111129002
$ws.Range("N132").Value = -17060

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1816.07
$ws.Range("J138").Value = 1881.2903
$ws.Range("L138").Value = 5643.8709
$ws.Range("N138").Value = -15923.8709

# Row 140: Tome for Tradition / Book of Ra'Kaznar
$ws.Range("H140").Value = 70780
$ws.Range("J140").Value = 70780
$ws.Range("L140").Value = 70780
$ws.Range("N140").Value = -81140

$ws = $wb.Worksheets.Item("ARM")
# Row 24: A Firm Hand / Iron Gauntlets
$ws.Range("H24").Value = 10691.25
$ws.Range("J24").Value = 10691.25
$ws.Range("L24").Value = 10691.25
$ws.Range("N24").Value = -11439.25

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 2797.2153
$ws.Range("I32").Value = 3069.5098
$ws.Range("K32").Value = 3069.5098
$ws.Range("M32").Value = -2782.5098

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 1754.6666
$ws.Range("I61").Value = 1625.4286
$ws.Range("J61").Value = 2207
$ws.Range("K61").Value = 1625.4286
$ws.Range("L61").Value = 2207
$ws.Range("M61").Value = -1413.4286
$ws.Range("N61").Value = -2631

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 828.52
$ws.Range("I74").Value = 715
$ws.Range("J74").Value = 1188
$ws.Range("K74").Value = 715
$ws.Range("L74").Value = 1188
$ws.Range("M74").Value = 159
$ws.Range("N74").Value = -2936

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 828.52
$ws.Range("I77").Value = 715
$ws.Range("J77").Value = 1188
$ws.Range("K77").Value = 3575
$ws.Range("L77").Value = 5940
$ws.Range("M77").Value = 793
$ws.Range("N77").Value = -14676

# Row 100: En Garde and on Guard / Doman Iron Gauntlets of Fending
$ws.Range("H100").Value = 10691.25
$ws.Range("J100").Value = 10691.25
$ws.Range("L100").Value = 10691.25
$ws.Range("N100").Value = -12855.25

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 55589188
$ws.Range("I102").Value = 55589188
$ws.Range("K102").Value = 55589188
$ws.Range("M102").Value = -55587566

# Row 114: A New Regular / Bluespirit Gauntlets of Fending
$ws.Range("H114").Value = 22599.2
$ws.Range("J114").Value = 22599.2
$ws.Range("L114").Value = 22599.2
$ws.Range("N114").Value = -31277.2

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 3196.889
$ws.Range("I132").Value = 2969.739
$ws.Range("J132").Value = 4503
$ws.Range("K132").Value = 8909.217000000001
$ws.Range("L132").Value = 13509
$ws.Range("M132").Value = -6379.217000000001
$ws.Range("N132").Value = -18569

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1754.6666
$ws.Range("I136").Value = 1625.4286
$ws.Range("J136").Value = 2207
$ws.Range("K136").Value = 4876.2858
$ws.Range("L136").Value = 6621
$ws.Range("M136").Value = -2326.2858
$ws.Range("N136").Value = -11721

# Row 139: Backing up My Words / Titanium Gold Thornplate of Fending
$ws.Range("H139").Value = 50715
$ws.Range("J139").Value = 50715
$ws.Range("L139").Value = 50715
$ws.Range("N139").Value = -60995

$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 26317066
$ws.Range("I99").Value = 31251278
$ws.Range("K99").Value = 31251278
$ws.Range("M99").Value = -31249780

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 1502.2413
$ws.Range("I107").Value = 1383.3
$ws.Range("J107").Value = 1766.5555
$ws.Range("K107").Value = 1383.3
$ws.Range("L107").Value = 1766.5555
$ws.Range("M107").Value = 536.7
$ws.Range("N107").Value = -5606.5555

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 8354.581
$ws.Range("I134").Value = 5842.5386
$ws.Range("J134").Value = 21417.2
$ws.Range("K134").Value = 17527.6158
$ws.Range("L134").Value = 64251.60000000001
$ws.Range("M134").Value = -14992.6158
$ws.Range("N134").Value = -69321.60000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 1443.5238
$ws.Range("I31").Value = 1180.3
$ws.Range("J31").Value = 1682.8182
$ws.Range("K31").Value = 1180.3
$ws.Range("L31").Value = 1682.8182
$ws.Range("M31").Value = -885.3
$ws.Range("N31").Value = -2272.8182

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 1443.5238
$ws.Range("I34").Value = 1180.3
$ws.Range("J34").Value = 1682.8182
$ws.Range("K34").Value = 1180.3
$ws.Range("L34").Value = 1682.8182
$ws.Range("M34").Value = -978.3
$ws.Range("N34").Value = -2086.8182

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 1553.3334
$ws.Range("I58").Value = 1425.7142
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 1425.7142
$ws.Range("L58").Value = 2000
$ws.Range("M58").Value = -1222.7142
$ws.Range("N58").Value = -2406

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 1097867
$ws.Range("I99").Value = 2025607.8
$ws.Range("J99").Value = 1446.1818
$ws.Range("K99").Value = 2025607.8
$ws.Range("L99").Value = 1446.1818
$ws.Range("M99").Value = -2024109.8
$ws.Range("N99").Value = -4442.1818

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 1097867
$ws.Range("I126").Value = 2025607.8
$ws.Range("J126").Value = 1446.1818
$ws.Range("K126").Value = 6076823.4
$ws.Range("L126").Value = 4338.5454
$ws.Range("M126").Value = -6074353.4
$ws.Range("N126").Value = -9278.545399999999

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 10131.733
$ws.Range("I132").Value = 13096.4
$ws.Range("J132").Value = 4202.4
$ws.Range("K132").Value = 39289.2
$ws.Range("L132").Value = 12607.2
$ws.Range("M132").Value = -36759.2
$ws.Range("N132").Value = -17667.2

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 16668384
$ws.Range("I134").Value = 18520094
$ws.Range("K134").Value = 55560282
$ws.Range("M134").Value = -55557747

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 1553.3334
$ws.Range("I136").Value = 1425.7142
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 4277.142599999999
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -1727.142599999999
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("CUL")
# Row 60: Drinking to Your Health / Mulled Tea
$ws.Range("H60").Value = 1889.8
$ws.Range("I60").Value = 771.5
$ws.Range("J60").Value = 2635.3333
$ws.Range("K60").Value = 2314.5
$ws.Range("L60").Value = 7905.999899999999
$ws.Range("M60").Value = -2063.5
$ws.Range("N60").Value = -8407.999899999999

# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 1781.0278
$ws.Range("I68").Value = 687.75
$ws.Range("J68").Value = 2093.3928
$ws.Range("K68").Value = 2063.25
$ws.Range("L68").Value = 6280.178400000001
$ws.Range("M68").Value = -1252.25
$ws.Range("N68").Value = -7902.178400000001

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 1781.0278
$ws.Range("I71").Value = 687.75
$ws.Range("J71").Value = 2093.3928
$ws.Range("K71").Value = 6189.75
$ws.Range("L71").Value = 18840.5352
$ws.Range("M71").Value = -2133.75
$ws.Range("N71").Value = -26952.5352

# Row 76: Old Victories, New Tastes / Dhalmel Fricassee
$ws.Range("H76").Value = 7070.615
$ws.Range("I76").Value = 10013
$ws.Range("J76").Value = 6825.4165
$ws.Range("K76").Value = 30039
$ws.Range("L76").Value = 20476.2495
$ws.Range("M76").Value = -29656
$ws.Range("N76").Value = -21242.2495

# Row 79: The Eats of Authenticity (L) / Dhalmel Fricassee
$ws.Range("H79").Value = 7070.615
$ws.Range("I79").Value = 10013
$ws.Range("J79").Value = 6825.4165
$ws.Range("K79").Value = 30039
$ws.Range("L79").Value = 20476.2495
$ws.Range("M79").Value = -28713
$ws.Range("N79").Value = -23128.2495

# Row 98: Sweet Kiss of Death / Rice Vinegar
$ws.Range("H98").Value = 1247.5
$ws.Range("J98").Value = 1597.5
$ws.Range("L98").Value = 4792.5
$ws.Range("N98").Value = -7788.5

# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 593.2857
$ws.Range("I122").Value = 499.5
$ws.Range("J122").Value = 663.625
$ws.Range("K122").Value = 4495.5
$ws.Range("L122").Value = 5972.625
$ws.Range("M122").Value = -2045.5
$ws.Range("N122").Value = -10872.625

$ws = $wb.Worksheets.Item("GSM")
# Row 94: Wants and Needles / Bombfish Needle
$ws.Range("H94").Value = 49711
$ws.Range("J94").Value = 49711
$ws.Range("L94").Value = 49711
$ws.Range("N94").Value = -51063

$ws = $wb.Worksheets.Item("LTW")
# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 5599.3184
$ws.Range("I136").Value = 10033.728
$ws.Range("J136").Value = 1164.909
$ws.Range("K136").Value = 30101.184
$ws.Range("L136").Value = 3494.727
$ws.Range("M136").Value = -27551.184
$ws.Range("N136").Value = -8594.727000000001

# Row 141: Just Generally Freezing / Gargantuaskin Trousers of Striking
$ws.Range("H141").Value = 49715
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 49715
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 49715
$ws.Range("M141").Value = ""
$ws.Range("N141").Value = -60075

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("H81").Value = 3237.1333
$ws.Range("I81").Value = 2696.6924
$ws.Range("J81").Value = 6750
$ws.Range("K81").Value = 5393.3848
$ws.Range("L81").Value = 13500
$ws.Range("M81").Value = -4332.3848
$ws.Range("N81").Value = -15622

# Row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("H84").Value = 3237.1333
$ws.Range("I84").Value = 2696.6924
$ws.Range("J84").Value = 6750
$ws.Range("K84").Value = 26966.924
$ws.Range("L84").Value = 67500
$ws.Range("M84").Value = -21662.924
$ws.Range("N84").Value = -78108

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 2392.75
$ws.Range("I132").Value = 2163
$ws.Range("K132").Value = 6489
$ws.Range("M132").Value = -3959

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 799.4
$ws.Range("I136").Value = 300
$ws.Range("J136").Value = 924.25
$ws.Range("K136").Value = 900
$ws.Range("L136").Value = 2772.75
$ws.Range("M136").Value = 1650
$ws.Range("N136").Value = -7872.75
